# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" list (column E, rows 16-38) was re-sorted in the
# underlying database from ascending (1608 .. 1807) to descending
# (1807 .. 1608). Re-apply the same reversal here so the displayed
# values end up in the new order while every other property of the
# cells (styles, merges, etc.) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 38

# Read the current "Periodo Mora" values into a plain array.
$periodos = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += $ws.Range("E$r").Value2
}

# Reverse the list manually (index based, since [array]::Reverse
# does not mutate the COM-backed array in place here).
$count = $periodos.Length
$reversed = @()
for ($i = $count - 1; $i -ge 0; $i--) {
    $reversed += $periodos[$i]
}

# Write the reversed order back into the same cells.
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $ws.Range("E$r").Value2 = $reversed[$i]
}
